# ThesisProposal.docx edit script
# Implements the changes described in the commit "Update word count and docs file"

$d = $word.ActiveDocument

$wdReplaceAll = 2
$LSQUO = [char]0x2018
$RSQUO = [char]0x2019
$LDQUO = [char]0x201C
$RDQUO = [char]0x201D
$ENDASH = [char]0x2013

function Find-ParagraphIndex($startsWith) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text.StartsWith($startsWith)) {
            return $i
        }
    }
    return -1
}

function Replace-Text($old, $new) {
    $d.Content.Find.Execute(
        $old, $false, $false, $false, $false, $false,
        $true, 1, $false, $new, $wdReplaceAll) | Out-Null
}

# ---------------------------------------------------------------------
# 1. Date: 2019-09-26 -> 2019-09-27
# ---------------------------------------------------------------------
Replace-Text "2019-09-26" "2019-09-27"

# ---------------------------------------------------------------------
# 2. Word count: 736 -> Word count:
# ---------------------------------------------------------------------
Replace-Text "Word count: 736" "Word count:"

# ---------------------------------------------------------------------
# 3. "resulting completed data sets" -> "completed data sets"
# ---------------------------------------------------------------------
Replace-Text `
    "The variability between the resulting completed data sets" `
    "The variability between the completed data sets"

# ---------------------------------------------------------------------
# 4. "a MI evaluation suite" -> "an MI evaluation suite"
# ---------------------------------------------------------------------
Replace-Text `
    "The goal is to develop a MI evaluation suite" `
    "The goal is to develop an MI evaluation suite"

Write-Host "Stage 1 (simple replacements) done"

# ---------------------------------------------------------------------
# 5. Rewrite the "Elaborate: ..." paragraph into the new discussion
#    paragraph, with two bold "source?" spans.
# ---------------------------------------------------------------------
$idxElaborate = Find-ParagraphIndex("Elaborate:")
$pElaborate = $d.Paragraphs($idxElaborate)
$fullRange = $pElaborate.Range.Duplicate
$fullRange.End = $fullRange.End - 1   # exclude the paragraph mark

$newParaText = (
    "The numerous assumptions underlying MI algorithms can principally not be verified from the incomplete observed data (" +
    "source?" +
    "). Instead, imputers are designated to evaluate the (plausibility of) multiply imputed data, and the robustness of results to the assumptions. The latter can be assessed by performing sensitivity analyses, see e.g. " +
    "Nguyen, Carlin, and Lee (2017)" +
    " " +
    "for practical guidelines. Methodology for evaluating the MI data however, is still largely missing. In the only available review of diagnostics," +
    " " +
    "Abayomi, Gelman, and Levy (2008)" +
    " " +
    "overlook the most vital state to be evaluated: convergence of the algorithm. Without convergence, any" +
    " " +
    $LSQUO + "deeper" + $RSQUO +
    " " +
    "assumption and resulting inference is invalid (" +
    "source?" +
    ")." +
    " "
)

$fullRange.Text = $newParaText
$fullRange.Font.Bold = 0

# Re-bold the two "source?" occurrences.
$parRange = $d.Paragraphs($idxElaborate).Range.Duplicate
$parRange.End = $parRange.End - 1

$search1 = $parRange.Duplicate
$search1.Find.Execute("source?") | Out-Null
$search1.Font.Bold = 1

$parEndNow = $d.Paragraphs($idxElaborate).Range.End - 1
$search2 = $d.Range($search1.End, $parEndNow)
$search2.Find.Execute("source?") | Out-Null
$search2.Font.Bold = 1

Write-Host "Stage 2 (Elaborate paragraph rewrite) done"

# ---------------------------------------------------------------------
# 6. Small text edits in the "Convergence properties..." paragraph.
# ---------------------------------------------------------------------
Replace-Text `
    "Convergence properties of iterative MI algorithms are still under debate" `
    "While the convergence properties of iterative MI algorithms are still under debate"

Replace-Text `
    ($ENDASH + "with specific procedures like") `
    ", specific procedures like"

Replace-Text `
    "posing an open question entirely" `
    "pose entirely open questions"

Replace-Text `
    "summarizes the state of the art as follows:" `
    "summarizes the issues with diagnosing convergence as follows:"

Write-Host "Stage 3 (convergence paragraph edits) done"

# ---------------------------------------------------------------------
# 7. Merge the "Van Buuren (2018) summarizes..." paragraph with the
#    following "Currently, applied researchers..." paragraph (which
#    becomes "Meanwhile, applied researchers...").
# ---------------------------------------------------------------------
$idxCur = Find-ParagraphIndex("Currently, applied researchers")
$idxPrev = $idxCur - 1

# Insert an extra leading space before "Currently" (becomes two spaces
# total once merged with the trailing space of the previous paragraph).
$pCur = $d.Paragraphs($idxCur)
$curStart = $pCur.Range.Start
$d.Range($curStart, $curStart).InsertBefore(" ") | Out-Null

Replace-Text `
    "Currently, applied researchers have to rely" `
    "Meanwhile, applied researchers have to rely"

# Delete the paragraph mark that separates the two paragraphs so they
# become one BodyText paragraph.
$prevEnd = $d.Paragraphs($idxPrev).Range.End
$d.Range($prevEnd - 1, $prevEnd).Delete() | Out-Null

Write-Host "Stage 4 (merge paragraphs, Currently -> Meanwhile) done"

# ---------------------------------------------------------------------
# 8. Simplify the "Add: which other assumptions..." paragraph.
# ---------------------------------------------------------------------
$idxAdd = Find-ParagraphIndex("Add: which other")
$pAdd = $d.Paragraphs($idxAdd)
$addRange = $pAdd.Range.Duplicate
$addRange.End = $addRange.End - 1
$addRange.Text = "Add: which other assumptions could be checked? The assumption of MAR. Ideally, we would want to study all possible combinations of variables: univariate, bivariate, etc. And to include both plots and stats."
$addRange.Font.Bold = 1

Write-Host "Stage 5 (Add: paragraph simplification) done"

# ---------------------------------------------------------------------
# 9. Insert a new bibliography entry "Nguyen, Carlin, and Lee (2017)"
#    right before the "Rubin, Donald B. 1987." reference.
# ---------------------------------------------------------------------
$idxRubin = Find-ParagraphIndex("Rubin, Donald B. 1987")
$pRubin = $d.Paragraphs($idxRubin)
$insertPoint = $d.Range($pRubin.Range.Start, $pRubin.Range.Start)
$insertPoint.InsertParagraphBefore() | Out-Null

$idxNguyen = $idxRubin
$pNguyen = $d.Paragraphs($idxNguyen)
$r = $pNguyen.Range.Duplicate
$r.End = $r.End - 1
$r.Text = "Nguyen, Cattram D., John B. Carlin, and Katherine J. Lee. 2017. " + $LDQUO + "Model Checking in Multiple Imputation: An Overview and Case Study." + $RDQUO
$r.Font.Italic = 0
$r.Font.Bold = 0

$parEnd = $d.Paragraphs($idxNguyen).Range.End - 1
$cursor = $d.Range($parEnd, $parEnd)
$cursor.InsertAfter(" ") | Out-Null

$parEnd = $d.Paragraphs($idxNguyen).Range.End - 1
$cursor = $d.Range($parEnd, $parEnd)
$italStart = $cursor.Start
$cursor.InsertAfter("Emerging Themes in Epidemiology") | Out-Null
$parEnd = $d.Paragraphs($idxNguyen).Range.End - 1
$italRange = $d.Range($italStart, $parEnd)
$italRange.Font.Italic = 1

$parEnd = $d.Paragraphs($idxNguyen).Range.End - 1
$cursor = $d.Range($parEnd, $parEnd)
$cursor.InsertAfter(" 14 (1): 8. ") | Out-Null

$parEnd = $d.Paragraphs($idxNguyen).Range.End - 1
$cursor = $d.Range($parEnd, $parEnd)
$hlStart = $cursor.Start
$cursor.InsertAfter("https://doi.org/10.1186/s12982-017-0062-6") | Out-Null
$parEnd = $d.Paragraphs($idxNguyen).Range.End - 1
$hlTextRange = $d.Range($hlStart, $parEnd)
$d.Hyperlinks.Add($hlTextRange, "https://doi.org/10.1186/s12982-017-0062-6") | Out-Null

$parEnd = $d.Paragraphs($idxNguyen).Range.End - 1
$cursor = $d.Range($parEnd, $parEnd)
$cursor.InsertAfter(".") | Out-Null

# Wrap the new paragraph's content in a "ref-nguy17" bookmark, mirroring
# the bookmarks around the other bibliography entries.
$bmRange = $d.Paragraphs($idxNguyen).Range.Duplicate
$bmRange.End = $bmRange.End - 1
$d.Bookmarks.Add("ref-nguy17", $bmRange) | Out-Null

Write-Host "Stage 6 (new Nguyen/Carlin/Lee reference) done"
Write-Host "New paragraph text: $($d.Paragraphs($idxNguyen).Range.Text)"
